# "Generate Report for Archive"
# The localization status for the a3b6c57b-...md file moved on from
# "Ready for handoff" to "In Translation" - update every place that
# status string appears (the Overview roll-up sheet as well as each
# per-locale detail sheet), then let the Status/locale columns re-fit
# to the new (shorter) text, same as Excel does automatically when the
# cell content driving a column's width changes.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: zh-cn (col E) and de-de (col F) status cells ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# --- Per-locale detail sheets: Status column (col C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus

# --- Re-fit the columns whose content just got shorter ---
# (Direct ColumnWidth assignment drives the same stored-width
# computation AutoFit uses, so nudge to the new natural width.)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
